$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the existing
# header style used by H1 ("IP").
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for the new I (I0) and J (IF) columns: row, I-value, J-value
$data = @(
    @(2, 4, 4),
    @(3, 7, 7),
    @(4, 9, 9),
    @(5, 5, 5),
    @(6, 7, 7),
    @(7, 4, 4),
    @(8, 6, 6),
    @(9, 6, 6),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 5, 5),
    @(13, 8, 8),
    @(14, 8, 8),
    @(15, 7, 8),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 7, 7),
    @(19, 6, 6),
    @(20, 6, 6),
    @(21, 7, 7),
    @(22, 8, 8),
    @(23, 6, 6),
    @(24, 6, 6),
    @(25, 7, 8),
    @(26, 5, 5),
    @(27, 9, 9),
    @(28, 7, 7),
    @(29, 8, 8),
    @(30, 5, 5),
    @(31, 7, 7),
    @(32, 10, 10),
    @(33, 7, 8),
    @(34, 7, 8),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 6, 6),
    @(38, 6, 7),
    @(39, 7, 7),
    @(40, 5, 5),
    @(41, 6, 6),
    @(42, 8, 8),
    @(43, 7, 8),
    @(44, 8, 8),
    @(45, 6, 6),
    @(46, 6, 6),
    @(47, 8, 8),
    @(48, 6, 6),
    @(49, 7, 7),
    @(50, 7, 7),
    @(51, 9, 9),
    @(52, 7, 7),
    @(53, 8, 8),
    @(54, 7, 8),
    @(55, 8, 8),
    @(56, 7, 7),
    @(57, 8, 8),
    @(58, 8, 8),
    @(59, 8, 8),
    @(60, 10, 10),
    @(61, 7, 7),
    @(62, 7, 7),
    @(63, 9, 9),
    @(64, 8, 8),
    @(65, 9, 9),
    @(66, 8, 8),
    @(67, 10, 10),
    @(68, 6, 6),
    @(69, 5, 6),
    @(70, 4, 4),
    @(71, 9, 9),
    @(72, 4, 4),
    @(73, 5, 5),
    @(74, 6, 6),
    @(75, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
